# Add three new market test-data sheets (Netherlands, Austria, Denmark)
# by cloning the existing "Greece" sheet (same layout/styles) and updating
# the sheet name plus the two market-specific cells (B4 = NGC ticket
# reference, B2 = market name).

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# --- Netherlands ---------------------------------------------------------
$source = $wb.Worksheets.Item("Greece")
$source.Copy($missing, $source)
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Value = "NGC-3144/T2188/T2189/T2190"
$netherlands.Range("B2").Value = "Netherlands Market"

# --- Austria --------------------------------------------------------------
$netherlands.Copy($missing, $netherlands)
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
$austria.Range("B4").Value = "NGC-3817/T2295"
$austria.Range("B2").Value = "Austria Market"

# --- Denmark ---------------------------------------------------------------
$austria.Copy($missing, $austria)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B4").Value = "NGC-2913/T2783"
$denmark.Range("B2").Value = "Denmark Market"

# Denmark ends up being the last/active tab, with the last used
# selection left on C21 (as captured in the source workbook).
$denmark.Activate()
$denmark.Range("C21").Select()

Write-Output "Added Netherlands, Austria and Denmark sheets"
